# "Added Key Thursday Speeches"
# Append the three new speech entries to the Thursday sheet of the DNC
# transcript metadata workbook, resize column C to fit the new (longer)
# filenames, and leave the workbook positioned on the Thursday tab.

$wb  = $excel.ActiveWorkbook
$thu = $wb.Worksheets.Item("Thursday")
$mon = $wb.Worksheets.Item("Monday")

# Row 2: Chelsea Clinton's speech
$thu.Cells.Item(2, 1).Value = "Clinton"
$thu.Cells.Item(2, 2).Value = "Chelsea"
$thu.Cells.Item(2, 3).Value = "dnc.clintonchelsea.txt"
$thu.Cells.Item(2, 4).Value = "Thursday"
$thu.Cells.Item(2, 5).Value = "speech"

# Row 3: Hillary Clinton's speech
$thu.Cells.Item(3, 1).Value = "Clinton"
$thu.Cells.Item(3, 2).Value = "Hillary"
$thu.Cells.Item(3, 3).Value = "dnc.clintonhillary.txt"
$thu.Cells.Item(3, 4).Value = "Thursday"
$thu.Cells.Item(3, 5).Value = "speech"

# Row 4: Khizr Khan's speech
$thu.Cells.Item(4, 1).Value = "Khan"
$thu.Cells.Item(4, 2).Value = "Khizer"
$thu.Cells.Item(4, 3).Value = "dnc.khan.txt"
$thu.Cells.Item(4, 4).Value = "Thursday"
$thu.Cells.Item(4, 5).Value = "speech"

# Widen column C so the longer filenames fit, matching the other day sheets.
$thu.Columns.Item(3).ColumnWidth = 25.17

# Move the Monday sheet's remembered selection down to B12 ...
$mon.Range("B12").Select() | Out-Null

# ... then switch to the Thursday sheet and leave the selection on A5,
# making Thursday the active tab when the workbook is saved.
$thu.Activate() | Out-Null
$thu.Range("A5").Select() | Out-Null
